# Auto-generated Excel COM-interop script to apply the 141 schedule update
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 07:49:14"
$ws1.Range("A3").Value = "Total filas: 98"

$ws1.Cells.Item(70, 1).Value = "07:49:14"
$ws1.Cells.Item(70, 2).Value = "07:49"
$ws1.Cells.Item(70, 3).Value = "10_OLMOS"
$ws1.Cells.Item(70, 4).Value = 0
$ws1.Cells.Item(70, 5).Value = "LP1912"

$ws1.Cells.Item(71, 1).Value = "07:49:14"
$ws1.Cells.Item(71, 2).Value = "07:51"
$ws1.Cells.Item(71, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(71, 4).Value = 2
$ws1.Cells.Item(71, 5).Value = "LP1912"

$ws1.Cells.Item(72, 1).Value = "07:18:13"
$ws1.Cells.Item(72, 2).Value = "07:52"
$ws1.Cells.Item(72, 3).Value = "215D_EL PATO"
$ws1.Cells.Item(72, 4).Value = 34
$ws1.Cells.Item(72, 5).Value = "LP1912"

$ws1.Cells.Item(73, 1).Value = "06:52:34"
$ws1.Cells.Item(73, 2).Value = "07:58"
$ws1.Cells.Item(73, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(73, 4).Value = 66
$ws1.Cells.Item(73, 5).Value = "LP1912"

$ws1.Cells.Item(74, 1).Value = "07:49:14"
$ws1.Cells.Item(74, 2).Value = "07:59"
$ws1.Cells.Item(74, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(74, 4).Value = 10
$ws1.Cells.Item(74, 5).Value = "LP1912"

$ws1.Cells.Item(75, 1).Value = "07:49:14"
$ws1.Cells.Item(75, 2).Value = "08:00"
$ws1.Cells.Item(75, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(75, 4).Value = 11
$ws1.Cells.Item(75, 5).Value = "LP1912"

$ws1.Cells.Item(76, 1).Value = "06:24:49"
$ws1.Cells.Item(76, 2).Value = "08:05"
$ws1.Cells.Item(76, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(76, 4).Value = 101
$ws1.Cells.Item(76, 5).Value = "LP1912"

$ws1.Cells.Item(77, 1).Value = "07:18:13"
$ws1.Cells.Item(77, 2).Value = "08:05"
$ws1.Cells.Item(77, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(77, 4).Value = 47
$ws1.Cells.Item(77, 5).Value = "LP1912"

$ws1.Cells.Item(78, 1).Value = "06:52:34"
$ws1.Cells.Item(78, 2).Value = "08:06"
$ws1.Cells.Item(78, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(78, 4).Value = 74
$ws1.Cells.Item(78, 5).Value = "LP1912"

$ws1.Cells.Item(79, 1).Value = "07:49:14"
$ws1.Cells.Item(79, 2).Value = "08:08"
$ws1.Cells.Item(79, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(79, 4).Value = 19
$ws1.Cells.Item(79, 5).Value = "LP1912"

$ws1.Cells.Item(80, 1).Value = "07:49:14"
$ws1.Cells.Item(80, 2).Value = "08:11"
$ws1.Cells.Item(80, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(80, 4).Value = 22
$ws1.Cells.Item(80, 5).Value = "LP1912"

$ws1.Cells.Item(81, 1).Value = "07:49:14"
$ws1.Cells.Item(81, 2).Value = "08:12"
$ws1.Cells.Item(81, 3).Value = "15_ABASTO"
$ws1.Cells.Item(81, 4).Value = 23
$ws1.Cells.Item(81, 5).Value = "LP1912"

$ws1.Cells.Item(82, 1).Value = "07:49:14"
$ws1.Cells.Item(82, 2).Value = "08:13"
$ws1.Cells.Item(82, 3).Value = "10_OLMOS"
$ws1.Cells.Item(82, 4).Value = 24
$ws1.Cells.Item(82, 5).Value = "LP1912"

$ws1.Cells.Item(83, 1).Value = "07:49:14"
$ws1.Cells.Item(83, 2).Value = "08:21"
$ws1.Cells.Item(83, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(83, 4).Value = 32
$ws1.Cells.Item(83, 5).Value = "LP1912"

$ws1.Cells.Item(84, 1).Value = "06:52:34"
$ws1.Cells.Item(84, 2).Value = "08:22"
$ws1.Cells.Item(84, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(84, 4).Value = 90
$ws1.Cells.Item(84, 5).Value = "LP1912"

$ws1.Cells.Item(85, 1).Value = "07:18:13"
$ws1.Cells.Item(85, 2).Value = "08:23"
$ws1.Cells.Item(85, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(85, 4).Value = 65
$ws1.Cells.Item(85, 5).Value = "LP1912"

$ws1.Cells.Item(86, 1).Value = "07:49:14"
$ws1.Cells.Item(86, 2).Value = "08:23"
$ws1.Cells.Item(86, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(86, 4).Value = 34
$ws1.Cells.Item(86, 5).Value = "LP1912"

$ws1.Cells.Item(87, 1).Value = "07:49:14"
$ws1.Cells.Item(87, 2).Value = "08:23"
$ws1.Cells.Item(87, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(87, 4).Value = 34
$ws1.Cells.Item(87, 5).Value = "LP1912"

$ws1.Cells.Item(88, 1).Value = "07:49:14"
$ws1.Cells.Item(88, 2).Value = "08:27"
$ws1.Cells.Item(88, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(88, 4).Value = 38
$ws1.Cells.Item(88, 5).Value = "LP1912"

$ws1.Cells.Item(89, 1).Value = "07:49:14"
$ws1.Cells.Item(89, 2).Value = "08:31"
$ws1.Cells.Item(89, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(89, 4).Value = 42
$ws1.Cells.Item(89, 5).Value = "LP1912"

$ws1.Cells.Item(90, 1).Value = "07:49:14"
$ws1.Cells.Item(90, 2).Value = "08:42"
$ws1.Cells.Item(90, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(90, 4).Value = 53
$ws1.Cells.Item(90, 5).Value = "LP1912"

$ws1.Cells.Item(91, 1).Value = "07:49:14"
$ws1.Cells.Item(91, 2).Value = "08:44"
$ws1.Cells.Item(91, 3).Value = "14_ABASTO"
$ws1.Cells.Item(91, 4).Value = 55
$ws1.Cells.Item(91, 5).Value = "LP1912"

$ws1.Cells.Item(92, 1).Value = "07:49:14"
$ws1.Cells.Item(92, 2).Value = "08:54"
$ws1.Cells.Item(92, 3).Value = "17_ROMERO"
$ws1.Cells.Item(92, 4).Value = 65
$ws1.Cells.Item(92, 5).Value = "LP1912"

$ws1.Cells.Item(93, 1).Value = "07:49:14"
$ws1.Cells.Item(93, 2).Value = "09:02"
$ws1.Cells.Item(93, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(93, 4).Value = 73
$ws1.Cells.Item(93, 5).Value = "LP1912"

$ws1.Cells.Item(94, 1).Value = "07:49:14"
$ws1.Cells.Item(94, 2).Value = "09:04"
$ws1.Cells.Item(94, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(94, 4).Value = 75
$ws1.Cells.Item(94, 5).Value = "LP1912"

$ws1.Cells.Item(95, 1).Value = "07:49:14"
$ws1.Cells.Item(95, 2).Value = "09:11"
$ws1.Cells.Item(95, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(95, 4).Value = 82
$ws1.Cells.Item(95, 5).Value = "LP1912"

$ws1.Cells.Item(96, 1).Value = "07:49:14"
$ws1.Cells.Item(96, 2).Value = "09:17"
$ws1.Cells.Item(96, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(96, 4).Value = 88
$ws1.Cells.Item(96, 5).Value = "LP1912"

$ws1.Cells.Item(97, 1).Value = "07:49:14"
$ws1.Cells.Item(97, 2).Value = "09:21"
$ws1.Cells.Item(97, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(97, 4).Value = 92
$ws1.Cells.Item(97, 5).Value = "LP1912"

$ws1.Cells.Item(98, 1).Value = "07:49:14"
$ws1.Cells.Item(98, 2).Value = "09:23"
$ws1.Cells.Item(98, 3).Value = "17_ROMERO"
$ws1.Cells.Item(98, 4).Value = 94
$ws1.Cells.Item(98, 5).Value = "LP1912"

$ws1.Cells.Item(99, 1).Value = "07:49:14"
$ws1.Cells.Item(99, 2).Value = "09:24"
$ws1.Cells.Item(99, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(99, 4).Value = 95
$ws1.Cells.Item(99, 5).Value = "LP1912"

$ws1.Cells.Item(100, 1).Value = "07:49:14"
$ws1.Cells.Item(100, 2).Value = "09:32"
$ws1.Cells.Item(100, 3).Value = "15_ABASTO"
$ws1.Cells.Item(100, 4).Value = 103
$ws1.Cells.Item(100, 5).Value = "LP1912"

$ws1.Cells.Item(101, 1).Value = "07:49:14"
$ws1.Cells.Item(101, 2).Value = "09:33"
$ws1.Cells.Item(101, 3).Value = "10_OLMOS"
$ws1.Cells.Item(101, 4).Value = 104
$ws1.Cells.Item(101, 5).Value = "LP1912"

$ws1.Cells.Item(102, 1).Value = "07:49:14"
$ws1.Cells.Item(102, 2).Value = "09:42"
$ws1.Cells.Item(102, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(102, 4).Value = 113
$ws1.Cells.Item(102, 5).Value = "LP1912"

$ws1.Cells.Item(103, 1).Value = "07:49:14"
$ws1.Cells.Item(103, 2).Value = "09:44"
$ws1.Cells.Item(103, 3).Value = "14_ABASTO"
$ws1.Cells.Item(103, 4).Value = 115
$ws1.Cells.Item(103, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 07:49:14"
$ws2.Range("A3").Value = "Total filas: 20"

$ws2.Cells.Item(21, 1).Value = "07:49:14"
$ws2.Cells.Item(21, 2).Value = "07:51"
$ws2.Cells.Item(21, 3).Value = "215D_EL PATO"
$ws2.Cells.Item(21, 4).Value = 2
$ws2.Cells.Item(21, 5).Value = "LP1912"

$ws2.Cells.Item(22, 1).Value = "07:18:13"
$ws2.Cells.Item(22, 2).Value = "07:52"
$ws2.Cells.Item(22, 3).Value = "215D_EL PATO"
$ws2.Cells.Item(22, 4).Value = 34
$ws2.Cells.Item(22, 5).Value = "LP1912"

$ws2.Cells.Item(23, 1).Value = "07:49:14"
$ws2.Cells.Item(23, 2).Value = "08:23"
$ws2.Cells.Item(23, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(23, 4).Value = 34
$ws2.Cells.Item(23, 5).Value = "LP1912"

$ws2.Cells.Item(24, 1).Value = "07:49:14"
$ws2.Cells.Item(24, 2).Value = "09:02"
$ws2.Cells.Item(24, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(24, 4).Value = 73
$ws2.Cells.Item(24, 5).Value = "LP1912"

$ws2.Cells.Item(25, 1).Value = "07:49:14"
$ws2.Cells.Item(25, 2).Value = "09:42"
$ws2.Cells.Item(25, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(25, 4).Value = 113
$ws2.Cells.Item(25, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 07:49:14"
$ws3.Range("A3").Value = "Total filas: 15"

$ws3.Cells.Item(17, 1).Value = "07:49:14"
$ws3.Cells.Item(17, 2).Value = "08:11"
$ws3.Cells.Item(17, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(17, 4).Value = 22
$ws3.Cells.Item(17, 5).Value = "L6203"

$ws3.Cells.Item(18, 1).Value = "06:52:34"
$ws3.Cells.Item(18, 2).Value = "08:33"
$ws3.Cells.Item(18, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18, 4).Value = 101
$ws3.Cells.Item(18, 5).Value = "L6173"

$ws3.Cells.Item(19, 1).Value = "07:49:14"
$ws3.Cells.Item(19, 2).Value = "08:35"
$ws3.Cells.Item(19, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(19, 4).Value = 46
$ws3.Cells.Item(19, 5).Value = "L6173"

$ws3.Cells.Item(20, 1).Value = "07:49:14"
$ws3.Cells.Item(20, 2).Value = "09:09"
$ws3.Cells.Item(20, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(20, 4).Value = 80
$ws3.Cells.Item(20, 5).Value = "L6203"
